$d = $word.ActiveDocument

# 1. Add new paragraph style "Abstract Title" (styleId AbstractTitle)
$abstractTitle = $d.Styles.Add("Abstract Title", 1)
$abstractTitle.BaseStyle = $d.Styles("Normal")
$abstractTitle.NextParagraphStyle = $d.Styles("Abstract")
$abstractTitle.QuickStyle = $true

$abstractTitle.ParagraphFormat.KeepWithNext = $true
$abstractTitle.ParagraphFormat.KeepTogether = $true
$abstractTitle.ParagraphFormat.Alignment = 1
$abstractTitle.ParagraphFormat.SpaceAfter = 0
$abstractTitle.ParagraphFormat.SpaceBefore = 15

$abstractTitle.Font.Size = 10
$abstractTitle.Font.SizeBi = 10
$abstractTitle.Font.Bold = $true
$abstractTitle.Font.Color = 9067060

# 2. Change the "Abstract" style's space-before from 300 twips (15pt) to 100 twips (5pt)
$abstract = $d.Styles("Abstract")
$abstract.ParagraphFormat.SpaceBefore = 5

Write-Output "Done"
